{"js": "// Circle Language Spec Plan: Set font to Calibri for non-heading text.\n// The built-in \"Normal\" paragraph style (used by all body/non-heading\n// text) is switched from Tahoma to Calibri at 11pt; heading styles are\n// left untouched.\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal,items/type\");\nawait context.sync();\n\nconst normalStyle = styles.items.find(\n  (s) => s.type === Word.StyleType.paragraph && s.nameLocal === \"Normal\"\n);\n\nif (!normalStyle) {\n  throw new Error('Could not find the built-in \"Normal\" style.');\n}\n\nnormalStyle.load(\"font\");\nawait context.sync();\n\nnormalStyle.font.name = \"Calibri\";\nnormalStyle.font.size = 11;\n\nawait context.sync();\n", "ps1": "# Circle Language Spec Plan: Set font to Calibri for non-heading text.\n# The built-in \"Normal\" paragraph style (used by all body/non-heading\n# text) is switched from Tahoma to Calibri at 11pt; heading styles are\n# left untouched.\n$d = $word.ActiveDocument\n\n$normalStyle = $d.Styles(\"Normal\")\n$normalStyle.Font.Name = \"Calibri\"\n$normalStyle.Font.Size = 11\n"}
